$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 9: SIAI blog post title/link updated
$ws.Range("D9").Value = "한국 인터넷에서 비난받는 SIAI 교육? (2)"
$ws.Range("E9").Value = "https://blog.pabii.co.kr/kr-internet-trollers-siai-2/#utm_source=rss&utm_medium=rss&utm_campaign=kr-internet-trollers-siai-2"

# Row 28: ropiens tistory post title/link updated
$ws.Range("D28").Value = "[논문리뷰] A fast, iterative solver for the Inverse Kinematcis Problem"
$ws.Range("E28").Value = "https://ropiens.tistory.com/199"

# Row 50: egloos post title/link updated
$ws.Range("D50").Value = "openCV super-resolution"
$ws.Range("E50").Value = "http://incredible.egloos.com/7558151"

# Row 51: bskyvision post title/link updated
$ws.Range("D51").Value = "[독후감] 금리와 환율의 변화가 세상을 어떻게 움직이는지 알 수 있게 해주는 <앞으로 3년 경제전쟁의 미래>"
$ws.Range("E51").Value = "https://bskyvision.com/entry/%EB%8F%85%ED%9B%84%EA%B0%90-%EA%B8%88%EB%A6%AC%EC%99%80-%ED%99%98%EC%9C%A8%EC%9D%98-%EB%B3%80%ED%99%94%EA%B0%80-%EC%84%B8%EC%83%81%EC%9D%84-%EC%96%B4%EB%96%BB%EA%B2%8C-%EC%9B%80%EC%A7%81%EC%9D%B4%EB%8A%94%EC%A7%80-%EC%95%8C-%EC%88%98-%EC%9E%88%EA%B2%8C-%ED%95%B4%EC%A3%BC%EB%8A%94-%EC%95%9E%EC%9C%BC%EB%A1%9C-3%EB%85%84-%EA%B2%BD%EC%A0%9C%EC%A0%84%EC%9F%81%EC%9D%98-%EB%AF%B8%EB%9E%98"
